$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BM (65): header "05-sep" plus one value per data row (2-18)
$ws.Cells.Item(1, 65).Value = "05-sep"

$ws.Cells.Item(2, 65).Value = 0
$ws.Cells.Item(3, 65).Value = 16.935527873406798
$ws.Cells.Item(4, 65).Value = 15.198662349629252
$ws.Cells.Item(5, 65).Value = 9.5836557024985272
$ws.Cells.Item(6, 65).Value = 0
$ws.Cells.Item(7, 65).Value = 15.149087125517745
$ws.Cells.Item(8, 65).Value = 11.134564455002007
$ws.Cells.Item(9, 65).Value = 10.710590896623469
$ws.Cells.Item(10, 65).Value = 15.316834092944312
$ws.Cells.Item(11, 65).Value = 12.470076000906772
$ws.Cells.Item(12, 65).Value = 0
$ws.Cells.Item(13, 65).Value = 9.3455170694668421
$ws.Cells.Item(14, 65).Value = 0
$ws.Cells.Item(15, 65).Value = 0
$ws.Cells.Item(16, 65).Value = 16.112116524939253
$ws.Cells.Item(17, 65).Value = 0
$ws.Cells.Item(18, 65).Value = 0

# Move selection to the newly added column, matching the author's saved view
$ws.Range("BM2:BM18").Select()
